$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New wallet-label rows (dates) appended below the existing Date/USDValue
# table. Format the cells as Text first so the date-look-alike strings are
# stored verbatim as shared strings (matching rows 2/3) instead of being
# auto-converted to date serials, then drop the temporary formatting so the
# cells end up with no explicit style, just like the existing rows.
$ws.Range("A4:A6").NumberFormat = "@"

$ws.Range("A4").Value = "2024-10-04"
$ws.Range("A5").Value = "2024-10-03"
$ws.Range("A6").Value = "2024-10-05"

$ws.Range("A4:A6").ClearFormats()
